# chore: update Sheets via scheduled runner
# Refreshes computed pricing/profit columns (H:N) on the Hades_Profits
# leve tables across the ALC/ARM/BSM/CRP/CUL/GSM/WVR sheets with newly
# scraped market-board data. A handful of rows flip between "has a
# computed LeveProfit" and "no price data" (LevePriceNQ/HQ == 0), so
# those LeveProfitNQ/HQ (M/N) cells are cleared instead of zeroed to
# match the existing convention (M/N is only populated when K/L != 0).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1748.7843
$ws.Range("I15").Value = 1748.7843
$ws.Range("K15").Value = 5246.3529
$ws.Range("M15").Value = -5077.3529
$ws.Range("H51").Value = 4506.25
$ws.Range("I51").Value = 500
$ws.Range("J51").Value = 5078.5713
$ws.Range("K51").Value = 500
$ws.Range("L51").Value = 5078.5713
$ws.Range("M51").Value = -16
$ws.Range("N51").Value = -6046.5713
$ws.Range("H74").Value = 3768.8
$ws.Range("I74").Value = 3537.6
$ws.Range("K74").Value = 3537.6
$ws.Range("M74").Value = -2601.6
$ws.Range("H77").Value = 3768.8
$ws.Range("I77").Value = 3537.6
$ws.Range("K77").Value = 17688
$ws.Range("M77").Value = -13008
$ws.Range("H113").Value = 2711.875
$ws.Range("I113").Value = 2775
$ws.Range("J113").Value = 2648.75
$ws.Range("K113").Value = 2775
$ws.Range("L113").Value = 2648.75
$ws.Range("M113").Value = 479
$ws.Range("N113").Value = -9156.75
$ws.Range("H138").Value = 2085305.2
$ws.Range("I138").Value = 1277.1052
$ws.Range("J138").Value = 2734428.8
$ws.Range("K138").Value = 3831.3156
$ws.Range("L138").Value = 8203286.399999999
$ws.Range("M138").Value = 1308.6844
$ws.Range("N138").Value = -8213566.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1368.625
$ws.Range("I2").Value = 1159.6842
$ws.Range("J2").Value = 2162.6
$ws.Range("K2").Value = 1159.6842
$ws.Range("L2").Value = 2162.6
$ws.Range("M2").Value = -1046.6842
$ws.Range("N2").Value = -2388.6
$ws.Range("H32").Value = 9839642
$ws.Range("I32").Value = 12387501
$ws.Range("J32").Value = 12187
$ws.Range("K32").Value = 12387501
$ws.Range("L32").Value = 12187
$ws.Range("M32").Value = -12387214
$ws.Range("N32").Value = -12761
$ws.Range("H61").Value = 167000900
$ws.Range("I61").Value = 250250900
$ws.Range("J61").Value = 500900
$ws.Range("K61").Value = 250250900
$ws.Range("L61").Value = 500900
$ws.Range("M61").Value = -250250688
$ws.Range("N61").Value = -501324
$ws.Range("H110").Value = 2503052.8
$ws.Range("I110").Value = 3337070.2
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 3337070.2
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = -3335025.2
$ws.Range("N110").Value = -5090
$ws.Range("H116").Value = 1368.625
$ws.Range("I116").Value = 1159.6842
$ws.Range("J116").Value = 2162.6
$ws.Range("K116").Value = 1159.6842
$ws.Range("L116").Value = 2162.6
$ws.Range("M116").Value = 1134.3158
$ws.Range("N116").Value = -6750.6
$ws.Range("H132").Value = 119441.35
$ws.Range("J132").Value = 335330
$ws.Range("L132").Value = 1005990
$ws.Range("N132").Value = -1011050
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 48011.72
$ws.Range("J134").Value = 48011.72
$ws.Range("L134").Value = 48011.72
$ws.Range("N134").Value = -58151.72
$ws.Range("H135").Value = 49602.07
$ws.Range("J135").Value = 49602.07
$ws.Range("L135").Value = 49602.07
$ws.Range("N135").Value = -59742.07
$ws.Range("H136").Value = 167000900
$ws.Range("I136").Value = 250250900
$ws.Range("J136").Value = 500900
$ws.Range("K136").Value = 750752700
$ws.Range("L136").Value = 1502700
$ws.Range("M136").Value = -750750150
$ws.Range("N136").Value = -1507800
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1368.625
$ws.Range("I3").Value = 1159.6842
$ws.Range("J3").Value = 2162.6
$ws.Range("K3").Value = 1159.6842
$ws.Range("L3").Value = 2162.6
$ws.Range("M3").Value = -1045.6842
$ws.Range("N3").Value = -2390.6
$ws.Range("H105").Value = 378150
$ws.Range("I105").Value = 603580
$ws.Range("J105").Value = 2433.3333
$ws.Range("K105").Value = 603580
$ws.Range("L105").Value = 2433.3333
$ws.Range("M105").Value = -601833
$ws.Range("N105").Value = -5927.3333
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1881.4667
$ws.Range("I16").Value = 1832.4615
$ws.Range("K16").Value = 1832.4615
$ws.Range("M16").Value = -1545.4615
$ws.Range("H96").Value = 20000
$ws.Range("J96").Value = 20000
$ws.Range("L96").Value = 20000
$ws.Range("N96").Value = -25492
$ws.Range("H107").Value = 988.8
$ws.Range("I107").Value = 444.84616
$ws.Range("J107").Value = 1999
$ws.Range("K107").Value = 444.84616
$ws.Range("L107").Value = 1999
$ws.Range("M107").Value = 1475.15384
$ws.Range("N107").Value = -5839
$ws.Range("H113").Value = 1881.4667
$ws.Range("I113").Value = 1832.4615
$ws.Range("K113").Value = 1832.4615
$ws.Range("M113").Value = 337.5385000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 880
$ws.Range("I92").Value = 300
$ws.Range("J92").Value = 1266.6666
$ws.Range("K92").Value = 900
$ws.Range("L92").Value = 3799.9998
$ws.Range("M92").Value = 348
$ws.Range("N92").Value = -6295.9998
$ws.Range("H124").Value = 1644.9259
$ws.Range("I124").Value = 749.5
$ws.Range("J124").Value = 1716.56
$ws.Range("K124").Value = 2248.5
$ws.Range("L124").Value = 5149.68
$ws.Range("M124").Value = 2661.5
$ws.Range("N124").Value = -14969.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 29918.9
$ws.Range("J115").Value = 29918.9
$ws.Range("L115").Value = 29918.9
$ws.Range("N115").Value = -33052.9
$ws.Range("H136").Value = 98067.14
$ws.Range("I136").Value = 70152.92999999999
$ws.Range("J136").Value = 167852.67
$ws.Range("K136").Value = 210458.79
$ws.Range("L136").Value = 503558.01
$ws.Range("M136").Value = -207908.79
$ws.Range("N136").Value = -508658.01
